$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "Synthese"
$ws2 = $wb.Worksheets.Item(2)   # "Par type de travaux"
$ws3 = $wb.Worksheets.Item(3)   # "Par entreprise"

# ---------------------------------------------------------------------------
# Sheet "Synthese": update the estimated/effective remaining budget figures,
# fix the "Depenses prevus" -> "Depenses prevues" typo and refresh the spend
# numbers now that the per-bill breakdown moved to its own sheet/model.
# ---------------------------------------------------------------------------
$ws1.Range("B2").Value = 492911
$ws1.Range("B3").Value = 492911
$ws1.Range("A4").Value = "Depenses prevues"
$ws1.Range("B4").Value = 0
$ws1.Range("B5").Value = 7089
$ws1.PageSetup.PaperSize = 9

# ---------------------------------------------------------------------------
# Sheet "Par type de travaux": drop the "espace vert" line, fix the
# "traveaux"/"prevus" typos and refresh electricite/plomberie totals.
# ---------------------------------------------------------------------------
$ws2.Rows.Item(3).Delete()

$ws2.Range("A1").Value = "Type de travaux"
$ws2.Range("B1").Value = "Depenses prevues"
$ws2.Range("C1").Value = "Depenses effectives"

$ws2.Range("A2").Value = "electricité"
$ws2.Range("B2").Value = 0
$ws2.Range("C2").Value = 5364

$ws2.Range("A3").Value = "plomberie"
$ws2.Range("B3").Value = 0
$ws2.Range("C3").Value = 1725

$ws2.PageSetup.PaperSize = 9

# ---------------------------------------------------------------------------
# Sheet "Par entreprise": bills now live per-company in their own rows;
# drop the old duplicate "je suis sous l'eau" row and refresh company names
# and totals to match the new per-type breakdown.
# ---------------------------------------------------------------------------
$ws3.Rows.Item(5).Delete()

$ws3.Range("A1").Value = "Nom de l'entreprise"
$ws3.Range("B1").Value = "Depenses prevues"
$ws3.Range("C1").Value = "Depenses effectives"

$ws3.Range("A2").Value = "electro 2000"
$ws3.Range("B2").Value = 0
$ws3.Range("C2").Value = 5364

$ws3.Range("A3").Value = "je suis sous l'eau"
$ws3.Range("B3").Value = 0
$ws3.Range("C3").Value = 725

$ws3.Range("A4").Value = "Robert & fils"
$ws3.Range("B4").Value = 0
$ws3.Range("C4").Value = 1000

$ws3.PageSetup.PaperSize = 9

# ---------------------------------------------------------------------------
# View state: restore each sheet's own selection, then make "Synthese" the
# active tab with G4 selected (matches the workbook's last-saved view).
# ---------------------------------------------------------------------------
$ws3.Activate() | Out-Null
$ws3.Range("A1").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("A1").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("G4").Select() | Out-Null
